$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 12.17857037013896
$ws.Range("C2").Value = 11.79194470087159
$ws.Range("D2").Value = 6.19782700849044
$ws.Range("F2").Value = 29.8353199738136
$ws.Range("G2").Value = 3.671543458454487
$ws.Range("K2").Value = 8.199793783997894
$ws.Range("L2").Value = 11.09203884864309
$ws.Range("M2").Value = 14.7303778700736
$ws.Range("O2").Value = 26.76033472215586
$ws.Range("B3").Value = 11.96517668011687
$ws.Range("C3").Value = 11.810476854906
$ws.Range("D3").Value = 6.161121354984201
$ws.Range("F3").Value = 29.87071967616454
$ws.Range("G3").Value = 3.673433943802937
$ws.Range("K3").Value = 8.033240175752855
$ws.Range("L3").Value = 11.10068201662378
$ws.Range("M3").Value = 14.70211086962703
$ws.Range("O3").Value = 26.82950733524058
$ws.Range("B4").Value = 11.83464196311016
$ws.Range("C4").Value = 11.82251492575145
$ws.Range("D4").Value = 6.138060997384338
$ws.Range("F4").Value = 29.89931537573925
$ws.Range("G4").Value = 3.674656707862278
$ws.Range("K4").Value = 7.93055796172069
$ws.Range("L4").Value = 11.10767087684619
$ws.Range("M4").Value = 14.68705483098625
$ws.Range("O4").Value = 26.87706400472029
$ws.Range("B5").Value = 11.78164319587799
$ws.Range("C5").Value = 11.82758682049067
$ws.Range("D5").Value = 6.128534788125751
$ws.Range("F5").Value = 29.9126912118404
$ws.Range("G5").Value = 3.675170632791104
$ws.Range("K5").Value = 7.888667287472749
$ws.Range("L5").Value = 11.11094225454307
$ws.Range("M5").Value = 14.68150247736357
$ws.Range("O5").Value = 26.89772032103919
$ws.Range("B6").Value = 11.77285669722892
$ws.Range("C6").Value = 11.82843906284682
$ws.Range("D6").Value = 6.126945256916471
$ws.Range("F6").Value = 29.91501623906639
$ws.Range("G6").Value = 3.67525691560151
$ws.Range("K6").Value = 7.881710273230976
$ws.Range("L6").Value = 11.11151104480833
$ws.Range("M6").Value = 14.68061586206824
$ws.Range("O6").Value = 26.90122733050747
$ws.Range("B7").Value = 11.83392631724652
$ws.Range("C7").Value = 11.82258265314724
$ws.Range("D7").Value = 6.13793304221524
$ws.Range("F7").Value = 29.89948879454776
$ws.Range("G7").Value = 3.674663575448859
$ws.Range("K7").Value = 7.929993121115989
$ws.Range("L7").Value = 11.10771328112603
$ws.Range("M7").Value = 14.68697758306782
$ws.Range("O7").Value = 26.8773374168349
$ws.Range("B8").Value = 12.10493350315489
$ws.Range("C8").Value = 11.79819806505938
$ws.Range("D8").Value = 6.185280967676862
$ws.Range("F8").Value = 29.84610084482558
$ws.Range("G8").Value = 3.67218245935829
$ws.Range("K8").Value = 8.142487187612875
$ws.Range("L8").Value = 11.09467021432461
$ws.Range("M8").Value = 14.72015700591191
$ws.Range("O8").Value = 26.78312904903049
$ws.Range("B9").Value = 12.6371521731457
$ws.Range("C9").Value = 11.75558792843673
$ws.Range("D9").Value = 6.273870619933921
$ws.Range("F9").Value = 29.79591164397029
$ws.Range("G9").Value = 3.667806707921017
$ws.Range("K9").Value = 8.553405924787262
$ws.Range("L9").Value = 11.08241637441881
$ws.Range("M9").Value = 14.80324667388724
$ws.Range("O9").Value = 26.63881130915112
$ws.Range("B10").Value = 13.02458900866846
$ws.Range("C10").Value = 11.72742551743562
$ws.Range("D10").Value = 6.336217121594922
$ws.Range("F10").Value = 29.79232621466007
$ws.Range("G10").Value = 3.664887275644726
$ws.Range("K10").Value = 8.848629994482513
$ws.Range("L10").Value = 11.08150235314773
$ws.Range("M10").Value = 14.8749569489224
$ws.Range("O10").Value = 26.55753102546429
$ws.Range("B11").Value = 13.19925033033208
$ws.Range("C11").Value = 11.71528956698494
$ws.Range("D11").Value = 6.363955204740038
$ws.Range("F11").Value = 29.79792186073095
$ws.Range("G11").Value = 3.663622639795459
$ws.Range("K11").Value = 8.980868059720926
$ws.Range("L11").Value = 11.08283368649401
$ws.Range("M11").Value = 14.90981836210712
$ws.Range("O11").Value = 26.52594856013663
$ws.Range("B12").Value = 13.26509775720957
$ws.Range("C12").Value = 11.71079059524586
$ws.Range("D12").Value = 6.374366667406283
$ws.Range("F12").Value = 29.80107847805353
$ws.Range("G12").Value = 3.663152826023445
$ws.Range("K12").Value = 9.030599251681403
$ws.Range("L12").Value = 11.08358806152331
$ws.Range("M12").Value = 14.92333442946568
$ws.Range("O12").Value = 26.51476587106353
$ws.Range("B13").Value = 13.25093036923801
$ws.Range("C13").Value = 11.71175523796794
$ws.Range("D13").Value = 6.372128521884206
$ws.Range("F13").Value = 29.8003525259917
$ws.Range("G13").Value = 3.66325360585657
$ws.Range("K13").Value = 9.019904805307185
$ws.Range("L13").Value = 11.08341448168197
$ws.Range("M13").Value = 14.9204096216116
$ws.Range("O13").Value = 26.51713968830543
$ws.Range("B14").Value = 13.20467387148596
$ws.Range("C14").Value = 11.71491749945503
$ws.Range("D14").Value = 6.364813627965623
$ws.Range("F14").Value = 29.79816077261387
$ws.Range("G14").Value = 3.663583806283438
$ws.Range("K14").Value = 8.984966643919829
$ws.Range("L14").Value = 11.08289074062171
$ws.Range("M14").Value = 14.91092407105041
$ws.Range("O14").Value = 26.5250129753259
$ws.Range("B15").Value = 13.17630039523334
$ws.Range("C15").Value = 11.71686704914813
$ws.Range("D15").Value = 6.360320939336578
$ws.Range("F15").Value = 29.79695333691399
$ws.Range("G15").Value = 3.663787244275062
$ws.Range("K15").Value = 8.963519771083647
$ws.Range("L15").Value = 11.08260248929699
$ws.Range("M15").Value = 14.90515467149834
$ws.Range("O15").Value = 26.52993680820021
$ws.Range("B16").Value = 13.01313694418921
$ws.Range("C16").Value = 11.72823217959072
$ws.Range("D16").Value = 6.334391594344551
$ws.Range("F16").Value = 29.79210583312491
$ws.Range("G16").Value = 3.664971195097554
$ws.Range("K16").Value = 8.839942271517232
$ws.Range("L16").Value = 11.08145042065089
$ws.Range("M16").Value = 14.87272313254315
$ws.Range("O16").Value = 26.55970365323023
$ws.Range("B17").Value = 12.91258808541612
$ws.Range("C17").Value = 11.73537694738457
$ws.Range("D17").Value = 6.318323139854271
$ws.Range("F17").Value = 29.79098213765705
$ws.Range("G17").Value = 3.665713724321537
$ws.Range("K17").Value = 8.763568430559168
$ws.Range("L17").Value = 11.08119046569031
$ws.Range("M17").Value = 14.85339607932518
$ws.Range("O17").Value = 26.57934687512518
$ws.Range("B18").Value = 12.85460949053912
$ws.Range("C18").Value = 11.73955001535769
$ws.Range("D18").Value = 6.309022445974133
$ws.Range("F18").Value = 29.79101616676543
$ws.Range("G18").Value = 3.666146780327689
$ws.Range("K18").Value = 8.719448930639274
$ws.Range("L18").Value = 11.08120546188742
$ws.Range("M18").Value = 14.84249097451695
$ws.Range("O18").Value = 26.59115261138107
$ws.Range("B19").Value = 12.83495596828904
$ws.Range("C19").Value = 11.74097387944472
$ws.Range("D19").Value = 6.305863417718256
$ws.Range("F19").Value = 29.79114457260755
$ws.Range("G19").Value = 3.666294432854036
$ws.Range("K19").Value = 8.70447948901967
$ws.Range("L19").Value = 11.08123882377872
$ws.Range("M19").Value = 14.8388352020059
$ws.Range("O19").Value = 26.59523694172487
$ws.Range("B20").Value = 12.923307227986
$ws.Range("C20").Value = 11.73460979699666
$ws.Range("D20").Value = 6.320039732758922
$ws.Range("F20").Value = 29.79103135189941
$ws.Range("G20").Value = 3.665634062976265
$ws.Range("K20").Value = 8.77171871215304
$ws.Range("L20").Value = 11.08120111853406
$ws.Range("M20").Value = 14.8554316596574
$ws.Range("O20").Value = 26.57720328751105
$ws.Range("B21").Value = 13.21826898108249
$ws.Range("C21").Value = 11.71398604715507
$ws.Range("D21").Value = 6.366964716979925
$ws.Range("F21").Value = 29.79877639768929
$ws.Range("G21").Value = 3.663486572437005
$ws.Range("K21").Value = 8.995238546878788
$ws.Range("L21").Value = 11.08303779355433
$ws.Range("M21").Value = 14.91370172345113
$ws.Range("O21").Value = 26.52267930408266
$ws.Range("B22").Value = 13.40930769594273
$ws.Range("C22").Value = 11.70107037414661
$ws.Range("D22").Value = 6.39709342470421
$ws.Range("F22").Value = 29.80988529625402
$ws.Range("G22").Value = 3.662135946882891
$ws.Range("K22").Value = 9.139293537665777
$ws.Range("L22").Value = 11.08569613731465
$ws.Range("M22").Value = 14.95361587339951
$ws.Range("O22").Value = 26.49157369174848
$ws.Range("B23").Value = 13.30752621405037
$ws.Range("C23").Value = 11.70791233301627
$ws.Range("D23").Value = 6.38106340048404
$ws.Range("F23").Value = 29.80340366540699
$ws.Range("G23").Value = 3.662851976932727
$ws.Range("K23").Value = 9.062609369284521
$ws.Range("L23").Value = 11.08414428765688
$ws.Range("M23").Value = 14.93214788346086
$ws.Range("O23").Value = 26.50776049535061
$ws.Range("B24").Value = 12.9184616335078
$ws.Range("C24").Value = 11.73495642179223
$ws.Range("D24").Value = 6.319263856273404
$ws.Range("F24").Value = 29.7910069836859
$ws.Range("G24").Value = 3.665670058660183
$ws.Range("K24").Value = 8.768034623872479
$ws.Range("L24").Value = 11.08119579007303
$ws.Range("M24").Value = 14.85451073071253
$ws.Range("O24").Value = 26.57817080688713
$ws.Range("B25").Value = 12.49351679657644
$ws.Range("C25").Value = 11.76656084581854
$ws.Range("D25").Value = 6.250375692077418
$ws.Range("F25").Value = 29.80364490478907
$ws.Range("G25").Value = 3.668938360444367
$ws.Range("K25").Value = 8.443195554783166
$ws.Range("L25").Value = 11.08430800810795
$ws.Range("M25").Value = 14.77887089765186
$ws.Range("O25").Value = 26.67351380495395
